$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Skills section updates ---

# Messaging: Kafka, JMS -> Kafka, JMS,
Replace-Text "Kafka, JMS" "Kafka, JMS,"

# Databases: PostgreSQL, MongoDB, MySQL, DynamoDB, Redis -> PostgreSQL, CockroachDB, Redis
Replace-Text "PostgreSQL, MongoDB, MySQL, DynamoDB, Redis" "PostgreSQL, CockroachDB, Redis"

# Cloud & Infrastructure: AWS (Lambda, S3, DynamoDB), Kubernetes, Docker -> AWS, GCP, Kubernetes, Docker
Replace-Text "AWS (Lambda, S3, DynamoDB), Kubernetes, Docker" "AWS, GCP, Kubernetes, Docker"

# API Technologies: REST, GraphQL, SWIFT messaging protocols -> REST, GraphQL, gRPC, SWIFT messaging protocols
Replace-Text "REST, GraphQL, SWIFT messaging protocols" "REST, GraphQL, gRPC, SWIFT messaging protocols"

# Remove the entire "Tools:" bullet paragraph from Skills
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Tools:*") {
        $p.Range.Delete()
        break
    }
}

# --- Interests section updates ---

# Climbing: -> Bouldering: (bold label run)
Replace-Text "Climbing:" "Bouldering:"

# Climbing description text
Replace-Text "Enthusiastic climber who values the sense of achievement and strong community bonds it brings. I enjoy taking climbing trips with friends whenever the opportunity arises, embracing the relaxed atmosphere and camaraderie it fosters." "Enthusiastic climber who values both the sense of achievement and the strong community bonds the sport fosters. I regularly take climbing trips, embracing the camaraderie and shared challenges."

# Cultural Experiences: -> Argentinian Tango: (bold label run)
Replace-Text "Cultural Experiences:" "Argentinian Tango:"

# Cultural Experiences description text
Replace-Text "I love learning new things such as Tango and have a keen interest in languages, particularly Spanish and Mandarin. Recently, I fulfilled a childhood dream of training Muay Thai in Northern Thailand. Engaging with different cultures through language, martial arts, and dance enriches my perspective and fuels my drive to grow." "A dance that cultivates connection, improvisation, and musicality. I appreciate how it develops deep partnership and encourages creative expression."

# Add a new "Languages:" bullet paragraph at the end of the Interests section
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$pStart = $newPara.Range.Start

$label = "Languages:"
$body = "Passionate about learning languages, particularly Spanish. Engaging with different cultures through language enriches my perspective of the world."

$r = $newPara.Range
$r.Text = $label
$labelEnd = $pStart + $label.Length

$spaceIns = $d.Range($labelEnd, $labelEnd)
$spaceIns.InsertAfter(" ")
$spaceEnd = $labelEnd + 1

$bodyIns = $d.Range($spaceEnd, $spaceEnd)
$bodyIns.InsertAfter($body)

# The label run inherited Bold formatting from the previous paragraph; make sure
# the space + body text after the label are not bold.
$afterLabel = $d.Range($labelEnd, $newPara.Range.End - 1)
$afterLabel.Font.Bold = 0
